$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.129080295562744
$ws.Range("B1").Value = 4.257708549499512
$ws.Range("C1").Value = 2.038357019424438
$ws.Range("D1").Value = 1.549264073371887
$ws.Range("E1").Value = 1.378419399261475
